$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Organization website (if available)" value (B10) from
# www.stat.kg to www.stat.gov.kg
$ws.Range("B10").Value = "www.stat.gov.kg"

# Update the active selection to B9
$ws.Range("B9").Select()
